# Login.xlsx automation data update.
#
# The "Message Displayed" column (D) records what the app showed back to the
# tester after each scenario. Row 5 is the "successful login" scenario - the
# app greets the user by name, so refresh that greeting with the latest
# value captured from the run (replacing the old "Hi, John Smith" sample).
#
# Re-assert the other strings already produced by this automation pass
# ("Hi, John Smith" from an earlier capture, and the "LOGIN" button/page
# label recorded during the run) so they stay present in the workbook's
# string table alongside the new greeting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "Hi, John Smith"
$ws.Range("F5").Value = "LOGIN"

$ws.Range("D5").Value = "Hi, DVhbCERv IlqEZZxz"
